# Natmi LR-pairs (Amh-Acvr1) re-run with a third cluster, "FAPs",
# added alongside the existing "ECs"/"sCs" clusters (per Dr Hou's advice).
# This expands the sending/target cluster cross-product from 2x2 to 3x3
# (rows 2-10) and updates the corresponding expression-specificity metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("ECs", "Amh", "Acvr1", "ECs", 3, 1, 0.9043393333333333, 2.713018, 0.5538994939397102, 0.55389949393971, 3, 1, 5.057757666666666, 15.173273, 0.173378811020062, 0.173378811020062, 4.573929196434888, 41.165362767914, 0.09603443568388097, 0.09603443568388097)
    ,@("ECs", "Amh", "Acvr1", "FAPs", 3, 1, 0.9043393333333333, 2.713018, 0.5538994939397102, 0.55389949393971, 3, 1, 14.247411, 42.742233, 0.4883980890531961, 0.4883980890531961, 12.884494165466, 115.960447489194, 0.2705234543676868, 0.2705234543676868)
    ,@("ECs", "Amh", "Acvr1", "sCs", 3, 1, 0.9043393333333333, 2.713018, 0.5538994939397102, 0.55389949393971, 3, 1, 9.866548666666667, 29.599646, 0.3382230999267418, 0.3382230999267418, 8.922708043514222, 80.304372391628, 0.1873416038881423, 0.1873416038881423)
    ,@("FAPs", "Amh", "Acvr1", "ECs", 1, 0.3333333333333333, 0.133128, 0.399384, 0.08153967112920636, 0.08153967112920636, 3, 1, 5.057757666666666, 15.173273, 0.173378811020062, 0.173378811020062, 0.6733291626479999, 6.059962463832, 0.01413725123134867, 0.01413725123134867)
    ,@("FAPs", "Amh", "Acvr1", "FAPs", 1, 0.3333333333333333, 0.133128, 0.399384, 0.08153967112920636, 0.08153967112920636, 3, 1, 14.247411, 42.742233, 0.4883980890531961, 0.4883980890531961, 1.896729331608, 17.070563984472, 0.03982381956153046, 0.03982381956153046)
    ,@("FAPs", "Amh", "Acvr1", "sCs", 1, 0.3333333333333333, 0.133128, 0.399384, 0.08153967112920636, 0.08153967112920636, 3, 1, 9.866548666666667, 29.599646, 0.3382230999267418, 0.3382230999267418, 1.313513890896, 11.821625018064, 0.02757860033632723, 0.02757860033632723)
    ,@("sCs", "Amh", "Acvr1", "ECs", 3, 1, 0.5952103333333333, 1.785631, 0.3645608349310836, 0.3645608349310835, 3, 1, 5.057757666666666, 15.173273, 0.173378811020062, 0.173378811020062, 3.010429626695888, 27.093866640263, 0.06320712410483235, 0.06320712410483235)
    ,@("sCs", "Amh", "Acvr1", "FAPs", 3, 1, 0.5952103333333333, 1.785631, 0.3645608349310836, 0.3645608349310835, 3, 1, 14.247411, 42.742233, 0.4883980890531961, 0.4883980890531961, 8.480206250446999, 76.32185625402299, 0.1780508151239789, 0.1780508151239789)
    ,@("sCs", "Amh", "Acvr1", "sCs", 3, 1, 0.5952103333333333, 1.785631, 0.3645608349310836, 0.3645608349310835, 3, 1, 9.866548666666667, 29.599646, 0.3382230999267418, 0.3382230999267418, 5.872671720736221, 52.854045486626, 0.1233028957022723, 0.1233028957022723)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
